$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped
# from 45172 (2023-09-03) to 45175 (2023-09-06) for every data row
# (rows 2 through 250).
for ($row = 2; $row -le 250; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
